# Insert a new data row at row 482, pushing the existing rows 482:512 down to 483:513.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(482).Insert()

$ws.Range("A482").Value = 10
$ws.Range("B482").Value = "Vega Modelo de Temuco"
$ws.Range("C482").Value = "La Araucanía"
$ws.Range("D482").Value = 44585
$ws.Range("E482").Value = 9
$ws.Range("F482").Value = 100112003
$ws.Range("G482").Value = "Ajo"
$ws.Range("H482").Value = "Chino"
$ws.Range("I482").Value = "Primera"
$ws.Range("J482").Value = 155
$ws.Range("K482").Value = 20000
$ws.Range("L482").Value = 20000
$ws.Range("M482").Value = 20000
$ws.Range("N482").Value = "`$/caja 10 kilos"
$ws.Range("O482").Value = "China"
$ws.Range("P482").Value = 2000
$ws.Range("Q482").Value = 10
$ws.Range("R482").Value = "Hortaliza"
